$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = '@'
    $c.Value = $val
    $c.Style = 'Normal'
}

$ws.Range('D2').Value = '72.323.59'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '4.059.00'
$ws.Range('E3').Value = '  +3.43%  '
$ws.Range('E4').Value = '  -0.11%  '
Set-TextValue 'D5' '523.04'
$ws.Range('E5').Value = '  -1.93%  '
Set-TextValue 'D6' '148.55'
$ws.Range('E6').Value = '  +2.80%  '
Set-TextValue 'D7' '0.625'
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('E8').Value = '  +0.14%  '
Set-TextValue 'D9' '0.741'
$ws.Range('E9').Value = '  +1.73%  '
Set-TextValue 'D10' '0.176'
$ws.Range('E10').Value = '  +1.83%  '
Set-TextValue 'D11' '0.0000340'
$ws.Range('E11').Value = '  +1.86%  '
Set-TextValue 'D12' '47.09'
$ws.Range('E12').Value = '  +10.90%  '
Set-TextValue 'D13' '10.79'
$ws.Range('E13').Value = '  +4.69%  '
$ws.Range('D14').Value = '4.697.89'
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('D15').Value = '4.041.83'
$ws.Range('E15').Value = '  +2.95%  '
Set-TextValue 'D16' '21.50'
$ws.Range('E16').Value = '  +8.28%  '
Set-TextValue 'D17' '14.35'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').Value = '72.147.56'
$ws.Range('E20').Value = '  +3.72%  '
Set-TextValue 'D21' '446.98'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D22' '3.55'
$ws.Range('E22').Value = '  +5.66%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D23' '96.12'
$ws.Range('E23').Value = '  +9.03%  '
Set-TextValue 'D24' '14.51'
$ws.Range('E24').Value = '  +0.51%  '
Set-TextValue 'D25' '12.16'
$ws.Range('E25').Value = '  +4.93%  '
Set-TextValue 'D26' '4.07'
$ws.Range('E26').Value = '  -1.74%  '
Set-TextValue 'D27' '11.33'
$ws.Range('E27').Value = '  +5.59%  '
Set-TextValue 'D28' '37.36'
$ws.Range('E28').Value = '  +2.30%  '
Set-TextValue 'D29' '5.78'
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D30' '709.95'
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '3.08'
$ws.Range('E31').Value = '  +8.45%  '
$ws.Range('E32').Value = '  +2.15%  '
Set-TextValue 'D33' '0.130'
$ws.Range('E33').Value = '  +3.20%  '
Set-TextValue 'D34' '6.93'
$ws.Range('E34').Value = '  +15.70%  '
Set-TextValue 'D35' '67.72'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').Value = '0.0₃0915'
$ws.Range('E36').Value = '  +8.56%  '
Set-TextValue 'D37' '0.450'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('B38').Value = 'ThetaToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D38' '3.69'
$ws.Range('E38').Value = '  +24.63%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D39' '40.90'
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('E40').Value = '  +3.45%  '
$ws.Range('E41').Value = '  -0.07%  '
Set-TextValue 'D42' '0.998'
$ws.Range('E42').Value = '  -0.26%  '
Set-TextValue 'D43' '0.0489'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D44' '3.11'
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D45' '2.82'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +4.47%  '
$ws.Range('E47').Value = '  +3.17%  '
Set-TextValue 'D48' '3.20'
$ws.Range('E48').Value = '  +1.66%  '
Set-TextValue 'D49' '9.23'
$ws.Range('E49').Value = '  +7.99%  '
Set-TextValue 'D50' '0.000277'
$ws.Range('E50').Value = '  +22.49%  '
Set-TextValue 'D51' '3.36'
$ws.Range('E51').Value = '  +1.23%  '
